$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Qty executed upto date" (column C) quantities
$ws.Range("C8").Value = 55
$ws.Range("C9").Value = 74
$ws.Range("C10").Value = 91
$ws.Range("C11").Value = 33
$ws.Range("C12").Value = 97
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 30
$ws.Range("C15").Value = 98
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 68

# Update corresponding "Upto date Amount" (column G) text values, keeping them as text
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "18944.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "42952.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21846.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "1088.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "690.00"

# Update Grand Total rows (G19/H19 and G21/H21)
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "85520.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "85520.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "85520.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "85520.00"
